$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) columns.
# D-column updates use NumberFormat "@" (Text) while assigning so that
# Excel does not auto-convert numeric-looking strings (e.g. "580.62")
# into numbers, then restore the "Normal" style so no stray number
# formatting is left behind on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.193.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.22%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.136.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.43%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.62"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.86%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.130.47"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.29%  "

$ws.Range("E9").Value = "  +0.71%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.10%  "

$ws.Range("E11").Value = "  +2.03%  "

$ws.Range("E12").Value = "  +0.20%  "

$ws.Range("E13").Value = "  +0.98%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.59%  "

$ws.Range("E15").Value = "  -0.64%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.655.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.43%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.185.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.24%  "

$ws.Range("E18").Value = "  -0.81%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.137.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.05%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "487.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.719"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.54%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.75%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.63%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.26%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.47%  "

$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.68%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.40"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.16%  "

$ws.Range("E31").Value = "  +2.24%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "29.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.33%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0000100"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.79%  "

$ws.Range("E34").Value = "  -2.90%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.96"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.27%  "

$ws.Range("E37").Value = "  -0.19%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.38"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.74%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.11"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.76%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.19"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.45%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.314"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.00%  "

$ws.Range("E42").Value = "  +2.19%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.67"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.82"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.81%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.861.59"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.56%  "

$ws.Range("E46").Value = "  +2.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0359"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "136.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.98%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.13"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.67%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.22%  "
